$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the total value in B2 from 87000 to 522000
$ws.Range("B2").Value = 522000
